$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-pulled data: update column F (dSF) values for the affected rows
$updates = @{
    4  = 0
    21 = -5
    23 = -1
    28 = -8
    29 = 8
    30 = 2
    32 = -1
    34 = 1
    36 = 2
    37 = -4
    45 = -3
    50 = -1
    53 = -5
    55 = 0
    59 = -3
    60 = -1
    67 = -4
    68 = 5
    69 = 2
    74 = 1
    75 = -10
    76 = 0
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
